$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append a new row (row 71) with the latest prediction/query data point,
# mirroring the pattern of the previous rows (A:C = 0, D = step value,
# E = predicted value, F = "query").
$row = 71
$ws.Cells.Item($row, 1).Value = 0
$ws.Cells.Item($row, 2).Value = 0
$ws.Cells.Item($row, 3).Value = 0
$ws.Cells.Item($row, 4).Value = 0.306122
$ws.Cells.Item($row, 5).Value = -20.61727907070679
$ws.Cells.Item($row, 6).Value = "query"
